$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Pace) entered first for the three new rows
$ws.Range("C6").Value = "ET"
$ws.Range("C7").Value = "ET"
$ws.Range("C8").Value = "ET"

# Column A (Reps) entered afterward, bottom-to-top
$ws.Range("A8").Value = "2.5, 4, 5.5"
$ws.Range("A7").Value = "2.5, 4, 4,5"
$ws.Range("A6").Value = "2.5, 4, 5"

# Remaining numeric columns
$ws.Range("B6").Value = 1
$ws.Range("D6").Value = 50

$ws.Range("B7").Value = 1
$ws.Range("D7").Value = 45

$ws.Range("B8").Value = 1
$ws.Range("D8").Value = 55

$ws.Range("A13").Select()
